$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.833.29"
$ws.Range("E2").Value = "  +0.90%  "

$ws.Range("D3").Value = "3.126.79"
$ws.Range("E3").Value = "  +1.01%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "533.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.63%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.24%  "

$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("D8").Value = "3.124.77"
$ws.Range("E8").Value = "  +1.00%  "

$ws.Range("E9").Value = "  +5.57%  "

$ws.Range("E10").Value = "  +0.38%  "

$ws.Range("E11").Value = "  +0.90%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.412"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.39%  "

$ws.Range("D13").Value = "3.665.96"
$ws.Range("E13").Value = "  +1.09%  "

$ws.Range("E14").Value = "  +2.04%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.93"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.78%  "

$ws.Range("E16").Value = "  +1.38%  "

$ws.Range("D17").Value = "57.941.21"
$ws.Range("E17").Value = "  +1.05%  "

$ws.Range("D18").Value = "3.123.73"
$ws.Range("E18").Value = "  +1.08%  "

$ws.Range("E19").Value = "  +3.11%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.72"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.04%  "

$ws.Range("E21").Value = "  +3.32%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "367.98"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.03%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.10%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.67"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.90%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.27"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.71%  "

$ws.Range("E26").Value = "  +1.73%  "

$ws.Range("E27").Value = "  +1.23%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.05%  "

$ws.Range("D29").Value = "0.0₃0865"
$ws.Range("E29").Value = "  -2.39%  "

$ws.Range("E30").Value = "  -0.41%  "

$ws.Range("E31").Value = "  +0.48%  "

$ws.Range("E32").Value = "  +1.63%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.42"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.90%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.16"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.93%  "

$ws.Range("E35").Value = "  +3.04%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "159.41"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.52%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.08"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.79%  "

$ws.Range("E38").Value = "  +5.44%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "25.44"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.96%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.67"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.31%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0672"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.48%  "

$ws.Range("D42").Value = "2.517.65"
$ws.Range("E42").Value = "  +6.28%  "

$ws.Range("E43").Value = "  +0.24%  "

$ws.Range("E44").Value = "  +0.61%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "37.82"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.80%  "

$ws.Range("E46").Value = "  +1.65%  "

$ws.Range("E47").Value = "  +0.06%  "

$ws.Range("E48").Value = "  +1.43%  "

$ws.Range("E49").Value = "  +3.37%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.72"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.54%  "

$ws.Range("E51").Value = "  -1.08%  "
